# "remote execution in docker"
#
# Adds a "browser" column to the testdata sheet (chrome/firefox), duplicates
# the two existing data rows so there are now two chrome/firefox pairs, and
# leaves behind the usual cursor/selection + outline-level bookkeeping that
# Excel records when a sheet like this is edited.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # "runmanager"
$ws2 = $wb.Worksheets.Item(2)   # "testdata"

# --- testdata: new header cell for the browser column ---------------------
$ws2.Range("E1").Value = "browser"

# --- testdata: fill in the browser value for the existing two rows --------
$ws2.Range("E2").Value = "chrome"
$ws2.Range("E3").Value = "firefox"

# --- testdata: append two more rows repeating the login / login2 pairing --
$ws2.Range("A4").Value = "login"
$ws2.Range("B4").Value = "Admin"
$ws2.Range("C4").Value = "admin123"
$ws2.Range("D4").Value = "yes"
$ws2.Range("E4").Value = "chrome"

$ws2.Range("A5").Value = "login2"
$ws2.Range("B5").Value = "Admin"
$ws2.Range("C5").Value = "admin123"
$ws2.Range("D5").Value = "yes"
$ws2.Range("E5").Value = "firefox"

# --- outline-level bookkeeping on the testdata sheet -----------------------
# (sheetFormatPr outlineLevelRow/outlineLevelCol go from 2/3 up to 4/4)
$ws2.Rows.Item(1).OutlineLevel = 4
$ws2.Columns.Item(1).OutlineLevel = 4

# --- cursor / selection bookkeeping ----------------------------------------
# Select on runmanager first so that testdata ends up as the active tab,
# matching the workbook's final state.
$ws1.Range("B5").Select() | Out-Null
$ws2.Range("C5").Select() | Out-Null
